# Update section numbering in column A for rows 15-24 on "Sheet 1"
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet 1")

$ws.Range("A15").Value = 8.1
$ws.Range("A16").Value = 8.2
$ws.Range("A17").Value = 8.3
$ws.Range("A18").Value = 9.1
$ws.Range("A19").Value = 9.2
$ws.Range("A20").Value = 9.3
$ws.Range("A21").Value = 9.4
$ws.Range("A22").Value = 9.5
$ws.Range("A23").Value = 9.6
$ws.Range("A24").Value = 9.7
